# Update the crypto price/volume table with freshly scraped values.
# The diff only touches the "Price" (column D) and "Volume(1h)" (column E)
# inline-string cells for rows 2-51; every other cell is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2'  = '28.933.38'
    'E2'  = '  -1.52%  '
    'D3'  = '1.912.17'
    'E3'  = '  -1.59%  '
    'D4'  = '1.003'
    'E4'  = '  +0.06%  '
    'D5'  = '324.88'
    'E5'  = '  -0.29%  '
    'E6'  = '  +0.00%  '
    'D7'  = '0.4591'
    'E7'  = '  -0.92%  '
    'D8'  = '0.3816'
    'E8'  = '  -1.41%  '
    'D9'  = '0.07717'
    'E9'  = '  -1.41%  '
    'D10' = '0.9812'
    'E10' = '  +0.58%  '
    'E11' = '  -2.35%  '
    'D12' = '1.913.99'
    'E12' = '  -1.23%  '
    'D13' = '6.950'
    'E13' = '  -1.82%  '
    'D14' = '5.666'
    'E14' = '  -1.59%  '
    'D15' = '0.07030'
    'E15' = '  -0.34%  '
    'E16' = '  -0.13%  '
    'D17' = '83.82'
    'E17' = '  -3.36%  '
    'D18' = '0.000009466'
    'E18' = '  -3.55%  '
    'D19' = '16.68'
    'E19' = '  -2.47%  '
    'D20' = '1.001'
    'E20' = '  -0.09%  '
    'D21' = '28.924.45'
    'E21' = '  -1.64%  '
    'E22' = '  -2.67%  '
    'E23' = '  -1.50%  '
    'D24' = '2.093'
    'E24' = '  -0.10%  '
    'D25' = '158.41'
    'E25' = '  +0.80%  '
    'D26' = '19.06'
    'E26' = '  -1.82%  '
    'D27' = '5.664'
    'E27' = '  -1.58%  '
    'D28' = '117.39'
    'E28' = '  -0.84%  '
    'D29' = '1.866'
    'E29' = '  +0.51%  '
    'D30' = '0.09289'
    'E30' = '  -0.72%  '
    'E31' = '  +1.02%  '
    'D32' = '5.081'
    'E33' = '  -4.05%  '
    'D34' = '3.156'
    'E34' = '  +1.06%  '
    'E35' = '  -0.72%  '
    'D36' = '1.159'
    'E36' = '  +0.52%  '
    'E37' = '  -0.01%  '
    'D38' = '0.02047'
    'E38' = '  -1.77%  '
    'D39' = '7.420'
    'E39' = '  -3.44%  '
    'D40' = '0.5495'
    'E40' = '  -3.06%  '
    'E41' = '  -1.39%  '
    'D42' = '2.858'
    'E42' = '  +4.70%  '
    'D43' = '9.332'
    'E43' = '  -0.84%  '
    'D44' = '0.5184'
    'E44' = '  -2.05%  '
    'D45' = '11.22'
    'E45' = '  -1.58%  '
    'D46' = '0.06906'
    'E46' = '  +0.50%  '
    'D47' = '2.098'
    'E47' = '  +0.44%  '
    'D48' = '0.000002578'
    'E48' = '  -8.04%  '
    'E49' = '  -1.98%  '
    'D50' = '110.64'
    'E50' = '  -0.65%  '
    'D51' = '0.2880'
    'E51' = '  -4.15%  '
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    # Force text storage so strings like "28.933.38" / "0.07717" are kept
    # as literal text (inline strings), matching the source workbook,
    # instead of being reinterpreted as numbers/dates by Excel.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
}
